$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.09"
$ws.Range("E2").Value = "'1.97%"
$ws.Range("D3").Value = "'36.44"
$ws.Range("E3").Value = "'-5.10%"
$ws.Range("D4").Value = "'5.042"
$ws.Range("E4").Value = "'1.08%"
$ws.Range("E5").Value = "'1.42%"
$ws.Range("D6").Value = "'2.108"
$ws.Range("E6").Value = "'-3.40%"
$ws.Range("D7").Value = "'7.933"
$ws.Range("E7").Value = "'-0.56%"
$ws.Range("D8").Value = "'4.061"
$ws.Range("E8").Value = "'1.36%"
$ws.Range("D9").Value = "'0.9221"
$ws.Range("E9").Value = "'0.56%"
$ws.Range("D10").Value = "'0.09601"
$ws.Range("E10").Value = "'5.78%"
$ws.Range("D11").Value = "'0.1876"
$ws.Range("E11").Value = "'4.98%"
$ws.Range("D12").Value = "'0.08711"
$ws.Range("E12").Value = "'2.49%"
$ws.Range("D13").Value = "'0.03501"
$ws.Range("E13").Value = "'-1.16%"
$ws.Range("D14").Value = "'0.09914"
$ws.Range("E14").Value = "'-0.17%"
$ws.Range("D15").Value = "'0.001431"
$ws.Range("E15").Value = "'-2.82%"
$ws.Range("D16").Value = "'0.005676"
$ws.Range("E16").Value = "'-0.24%"
$ws.Range("D17").Value = "'3.463"
$ws.Range("E17").Value = "'-0.51%"
$ws.Range("D18").Value = "'2.412"
$ws.Range("E18").Value = "'8.83%"
$ws.Range("D19").Value = "'0.3417"
$ws.Range("E19").Value = "'-1.34%"
$ws.Range("E20").Value = "'1.91%"
$ws.Range("D21").Value = "'4.781"
$ws.Range("E21").Value = "'5.01%"
$ws.Range("D23").Value = "'0.04604"
$ws.Range("E23").Value = "'-1.31%"
$ws.Range("E24").Value = "'15.20%"
$ws.Range("D25").Value = "'0.001231"
$ws.Range("E25").Value = "'-0.16%"
$ws.Range("D26").Value = "'0.0001401"
$ws.Range("E26").Value = "'7.48%"
$ws.Range("E27").Value = "'-42.77%"
$ws.Range("D39").Value = "'0.01836"
$ws.Range("E39").Value = "'5.74%"
$ws.Range("D40").Value = "'0.04773"
$ws.Range("E40").Value = "'1.97%"
$ws.Range("D41").Value = "'0.007486"
$ws.Range("E41").Value = "'-5.67%"
$ws.Range("D42").Value = "'0.1403"
$ws.Range("E42").Value = "'1.26%"
$ws.Range("D43").Value = "'0.007741"
$ws.Range("E43").Value = "'0.47%"
$ws.Range("D44").Value = "'0.002232"
$ws.Range("E44").Value = "'1.37%"
$ws.Range("E45").Value = "'8.63%"
$ws.Range("D46").Value = "'0.00006177"
$ws.Range("E46").Value = "'2.16%"
$ws.Range("E47").Value = "'-0.20%"
$ws.Range("D48").Value = "'0.0005799"
$ws.Range("E48").Value = "'-0.03%"
$ws.Range("D49").Value = "'38.56"
$ws.Range("E49").Value = "'357.80%"
$ws.Range("D50").Value = "'0.001999"
$ws.Range("E50").Value = "'-25.91%"
$ws.Range("D51").Value = "'0.00002101"
$ws.Range("E51").Value = "'-0.20%"
